# Updates cryptos list values per Oct 2 2023 GitHub Actions data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'28.268.80"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "'  +3.85%  "
$ws.Cells.Item(2, 5).Style = "Normal"

$ws.Cells.Item(3, 4).Value = "'1.729.77"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "'  +2.33%  "
$ws.Cells.Item(3, 5).Style = "Normal"

$ws.Cells.Item(4, 5).Value = "'  -0.06%  "
$ws.Cells.Item(4, 5).Style = "Normal"

$ws.Cells.Item(5, 4).Value = "'219.25"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "'  +1.31%  "
$ws.Cells.Item(5, 5).Style = "Normal"

$ws.Cells.Item(6, 4).Value = "'0.522"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "'  +0.14%  "
$ws.Cells.Item(6, 5).Style = "Normal"

$ws.Cells.Item(7, 5).Value = "'  +0.00%  "
$ws.Cells.Item(7, 5).Style = "Normal"

$ws.Cells.Item(8, 4).Value = "'23.94"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "'  +4.04%  "
$ws.Cells.Item(8, 5).Style = "Normal"

$ws.Cells.Item(9, 5).Value = "'  +1.61%  "
$ws.Cells.Item(9, 5).Style = "Normal"

$ws.Cells.Item(10, 5).Value = "'  +1.39%  "
$ws.Cells.Item(10, 5).Style = "Normal"

$ws.Cells.Item(11, 5).Value = "'  +0.67%  "
$ws.Cells.Item(11, 5).Style = "Normal"

$ws.Cells.Item(12, 4).Value = "'1.976.47"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "'  +2.37%  "
$ws.Cells.Item(12, 5).Style = "Normal"

$ws.Cells.Item(13, 4).Value = "'1.730.40"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "'  +2.07%  "
$ws.Cells.Item(13, 5).Style = "Normal"

$ws.Cells.Item(14, 4).Value = "'4.24"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "'  +0.88%  "
$ws.Cells.Item(14, 5).Style = "Normal"

$ws.Cells.Item(15, 5).Value = "'  +1.57%  "
$ws.Cells.Item(15, 5).Style = "Normal"

$ws.Cells.Item(16, 4).Value = "'67.68"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "'  +0.14%  "
$ws.Cells.Item(16, 5).Style = "Normal"

$ws.Cells.Item(17, 4).Value = "'28.268.83"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "'  +3.78%  "
$ws.Cells.Item(17, 5).Style = "Normal"

$ws.Cells.Item(18, 4).Value = "'245.45"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "'  +3.35%  "
$ws.Cells.Item(18, 5).Style = "Normal"

$ws.Cells.Item(19, 4).Value = "'0.0₃0752"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "'  +0.91%  "
$ws.Cells.Item(19, 5).Style = "Normal"

$ws.Cells.Item(20, 4).Value = "'7.89"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "'  -3.27%  "
$ws.Cells.Item(20, 5).Style = "Normal"

$ws.Cells.Item(21, 5).Value = "'  -0.07%  "
$ws.Cells.Item(21, 5).Style = "Normal"

$ws.Cells.Item(22, 4).Value = "'4.65"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "'  +1.47%  "
$ws.Cells.Item(22, 5).Style = "Normal"

$ws.Cells.Item(23, 4).Value = "'9.68"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "'  +0.35%  "
$ws.Cells.Item(23, 5).Style = "Normal"

$ws.Cells.Item(24, 5).Value = "'  -0.84%  "
$ws.Cells.Item(24, 5).Style = "Normal"

$ws.Cells.Item(25, 4).Value = "'149.56"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "'  +0.88%  "
$ws.Cells.Item(25, 5).Style = "Normal"

$ws.Cells.Item(26, 4).Value = "'7.49"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "'  +2.35%  "
$ws.Cells.Item(26, 5).Style = "Normal"

$ws.Cells.Item(27, 4).Value = "'16.63"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "'  +0.79%  "
$ws.Cells.Item(27, 5).Style = "Normal"

$ws.Cells.Item(28, 5).Value = "'  +0.19%  "
$ws.Cells.Item(28, 5).Style = "Normal"

$ws.Cells.Item(29, 5).Value = "'  +0.24%  "
$ws.Cells.Item(29, 5).Style = "Normal"

$ws.Cells.Item(30, 5).Value = "'  +2.56%  "
$ws.Cells.Item(30, 5).Style = "Normal"

$ws.Cells.Item(31, 5).Value = "'  +2.26%  "
$ws.Cells.Item(31, 5).Style = "Normal"

$ws.Cells.Item(32, 4).Value = "'3.42"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "'  +0.47%  "
$ws.Cells.Item(32, 5).Style = "Normal"

$ws.Cells.Item(33, 5).Value = "'  +0.37%  "
$ws.Cells.Item(33, 5).Style = "Normal"

$ws.Cells.Item(34, 4).Value = "'1.487.58"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "'  -5.46%  "
$ws.Cells.Item(34, 5).Style = "Normal"

$ws.Cells.Item(35, 4).Value = "'1.65"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "'  -2.11%  "
$ws.Cells.Item(35, 5).Style = "Normal"

$ws.Cells.Item(36, 4).Value = "'0.974"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "'  +2.09%  "
$ws.Cells.Item(36, 5).Style = "Normal"

$ws.Cells.Item(37, 4).Value = "'0.602"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "'  -0.48%  "
$ws.Cells.Item(37, 5).Style = "Normal"

$ws.Cells.Item(38, 4).Value = "'2.41"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "'  +0.70%  "
$ws.Cells.Item(38, 5).Style = "Normal"

$ws.Cells.Item(39, 4).Value = "'0.0176"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "'  +0.82%  "
$ws.Cells.Item(39, 5).Style = "Normal"

$ws.Cells.Item(40, 5).Value = "'  +0.56%  "
$ws.Cells.Item(40, 5).Style = "Normal"

$ws.Cells.Item(41, 4).Value = "'69.90"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "'  +0.45%  "
$ws.Cells.Item(41, 5).Style = "Normal"

$ws.Cells.Item(42, 5).Value = "'  +0.00%  "
$ws.Cells.Item(42, 5).Style = "Normal"

$ws.Cells.Item(43, 4).Value = "'5.65"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "'  -2.10%  "
$ws.Cells.Item(43, 5).Style = "Normal"

$ws.Cells.Item(44, 4).Value = "'2.29"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "'  +1.35%  "
$ws.Cells.Item(44, 5).Style = "Normal"

$ws.Cells.Item(45, 4).Value = "'1.878.83"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "'  +1.93%  "
$ws.Cells.Item(45, 5).Style = "Normal"

$ws.Cells.Item(46, 5).Value = "'  +1.54%  "
$ws.Cells.Item(46, 5).Style = "Normal"

$ws.Cells.Item(47, 4).Value = "'1.72"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "'  +6.82%  "
$ws.Cells.Item(47, 5).Style = "Normal"

$ws.Cells.Item(48, 2).Value = "'BabyDogeCoin"
$ws.Cells.Item(48, 2).Style = "Normal"
$ws.Cells.Item(48, 3).Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(48, 3).Style = "Normal"
$ws.Cells.Item(48, 4).Value = "'0.0₆0114"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "'  +3.49%  "
$ws.Cells.Item(48, 5).Style = "Normal"

$ws.Cells.Item(49, 2).Value = "'Quant"
$ws.Cells.Item(49, 2).Style = "Normal"
$ws.Cells.Item(49, 3).Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(49, 3).Style = "Normal"
$ws.Cells.Item(49, 4).Value = "'90.34"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "'  -0.83%  "
$ws.Cells.Item(49, 5).Style = "Normal"

$ws.Cells.Item(50, 4).Value = "'8.16"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "'  -0.68%  "
$ws.Cells.Item(50, 5).Style = "Normal"

$ws.Cells.Item(51, 5).Value = "'  -1.12%  "
$ws.Cells.Item(51, 5).Style = "Normal"
